$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime columns (E, H) for rows 2-3
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 10:17:42"
$wsZh.Range("E3").Value = "2016-03-19 10:17:42"
$wsZh.Range("H2").Value = "2016-03-19 10:18:01"
$wsZh.Range("H3").Value = "2016-03-19 10:18:01"

# de-de sheet: update Correspond Handoff/Handback Datetime columns (E, H) for rows 2-3
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 10:17:46"
$wsDe.Range("E3").Value = "2016-03-19 10:17:46"
$wsDe.Range("H2").Value = "2016-03-19 10:18:07"
$wsDe.Range("H3").Value = "2016-03-19 10:18:07"
